$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing data rows (2..7) into memory before we start overwriting them,
# since the new row 2 pushes all the current records down by one row.
$old = @{}
for ($r = 2; $r -le 7; $r++) {
    $rowvals = @()
    for ($c = 1; $c -le 20; $c++) {
        $rowvals += ,$ws.Cells.Item($r, $c).Value2
    }
    $old[$r] = $rowvals
}

# Write the shifted rows back out, from the bottom up, into rows 3..8.
for ($r = 7; $r -ge 2; $r--) {
    $rowvals = $old[$r]
    $dst = $r + 1
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($dst, $c).Value = $rowvals[$c - 1]
    }
}

# Make sure the date cell in the newly-created row 8 keeps the same date format
# used by column D elsewhere in the sheet.
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Now write the brand-new record into row 2.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(2, 4).Value = 44545
$ws.Cells.Item(2, 5).Value = 15
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100103
$ws.Cells.Item(2, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(2, 9).Value = 100103003
$ws.Cells.Item(2, 10).Value = "Damasco"
$ws.Cells.Item(2, 11).Value = "Castle Brite"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 24000
$ws.Cells.Item(2, 15).Value = 25000
$ws.Cells.Item(2, 16).Value = 24500
$ws.Cells.Item(2, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(2, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(2, 19).Value = 1361
$ws.Cells.Item(2, 20).Value = 18
